# Engineering Yield Tracker - BFT sheet updates
# - Fix L4 (day 6) value/format to a proper percentage
# - Fill in O4 (day 9) and P4 (day 10) FPY values
# - Clear out the stale "future" data that had been entered in columns AG:AK
#   (days 27-31) across the BFT table - these were placeholder/test values
#   that shouldn't have been there yet.
# - Move the active selection to P8

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("BFT")

# Fix the FPY value for day 6 (L4): was stored as "99.21" under a General
# format; correct it to the 0.9921 percentage value with a percent format.
$L4 = $ws.Range("L4")
$L4.Value = 0.9921
$L4.NumberFormat = "0.00%"

# Populate day 9 (O4) and day 10 (P4) FPY values.
$O4 = $ws.Range("O4")
$O4.Value = 0.987
$O4.NumberFormat = "0.00%"

$P4 = $ws.Range("P4")
$P4.Value = 0.9793
$P4.NumberFormat = "0.00%"

# Clear the (erroneous / not-yet-applicable) data that had been filled into
# columns AG:AK (days 27-31) for the Monica customer rows (4-9).
$ws.Range("AG4:AK9").ClearContents()

# Clear AG:AH (days 27-28) for the Fiona customer rows (10-20).
$ws.Range("AG10:AH20").ClearContents()

# Clear AG (day 27) for the first few AMBER/TETON rows (21-23).
$ws.Range("AG21:AG23").ClearContents()

# Update the active cell selection.
$ws.Range("P8").Select()
